# Atualização automática de SOLEDADE.xlsx
#
# Changes applied (per the authoritative diff):
#   1. Rename "Paineis DARQ"            -> "PAINEIS DARQ"
#   2. Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
#   3. Delete the "Desarquivamentos Pendentes" sheet entirely
#
# (Deleting sheet 3 above also removes the shared strings and cell
#  formats that were exclusively used by it, matching the rest of the
#  diff as a natural consequence.)

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$painelSheet = $wb.Worksheets.Item("Paineis DARQ")
$painelSheet.Name = "PAINEIS DARQ"

$recolhimentoSheet = $wb.Worksheets.Item("Recolhimento x Eliminacao")
$recolhimentoSheet.Name = "RECOLHIMENTO X ELIMINAÇÃO"

$desarquivamentosSheet = $wb.Worksheets.Item("Desarquivamentos Pendentes")
$desarquivamentosSheet.Delete()

$excel.DisplayAlerts = $true
